# Swap the "step"/"expected result" content between the TC2 test-case
# block (row 18) and the TC3 test-case block (row 25), so that the
# "analisar prestação de contas" content now appears under TC2 and the
# "cancelar diária" content now appears under TC3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tc2Step   = $ws.Range("B18").Value2
$tc2Result = $ws.Range("D18").Value2
$tc3Step   = $ws.Range("B25").Value2
$tc3Result = $ws.Range("D25").Value2

$ws.Range("B18").Value = $tc3Step
$ws.Range("D18").Value = $tc3Result
$ws.Range("B25").Value = $tc2Step
$ws.Range("D25").Value = $tc2Result
